$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.87"
$ws.Range("E2").Value = "'1.53%"
$ws.Range("D3").Value = "'27.45"
$ws.Range("E3").Value = "'1.66%"
$ws.Range("D4").Value = "'4.709"
$ws.Range("E4").Value = "'2.69%"
$ws.Range("D5").Value = "'0.06068"
$ws.Range("E5").Value = "'2.78%"
$ws.Range("D6").Value = "'6.709"
$ws.Range("E6").Value = "'1.15%"
$ws.Range("D7").Value = "'0.8632"
$ws.Range("E7").Value = "'0.93%"
$ws.Range("D8").Value = "'0.9243"
$ws.Range("E8").Value = "'-1.38%"
$ws.Range("D9").Value = "'0.1409"
$ws.Range("E9").Value = "'0.24%"
$ws.Range("D10").Value = "'0.05093"
$ws.Range("E10").Value = "'6.02%"
$ws.Range("D11").Value = "'0.07111"
$ws.Range("E11").Value = "'0.20%"
$ws.Range("D12").Value = "'0.03073"
$ws.Range("E12").Value = "'-1.64%"
$ws.Range("D13").Value = "'0.09105"
$ws.Range("E13").Value = "'-0.50%"
$ws.Range("D14").Value = "'0.001529"
$ws.Range("E14").Value = "'-0.52%"
$ws.Range("D15").Value = "'0.0006077"
$ws.Range("E15").Value = "'-94.21%"
$ws.Range("D16").Value = "'0.006197"
$ws.Range("E16").Value = "'-0.43%"
$ws.Range("D17").Value = "'3.473"
$ws.Range("E17").Value = "'-1.36%"
$ws.Range("D18").Value = "'3.167"
$ws.Range("E18").Value = "'-0.52%"
$ws.Range("E19").Value = "'-2.21%"
$ws.Range("E20").Value = "'2.38%"
$ws.Range("D21").Value = "'0.1307"
$ws.Range("E21").Value = "'2.27%"
$ws.Range("D22").Value = "'4.099"
$ws.Range("E22").Value = "'7.23%"
$ws.Range("D23").Value = "'0.04252"
$ws.Range("E23").Value = "'-0.69%"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("E24").Value = "'-0.43%"
$ws.Range("D25").Value = "'0.003914"
$ws.Range("E25").Value = "'-8.86%"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'0.05%"
$ws.Range("E40").Value = "'1.52%"
$ws.Range("D41").Value = "'0.1116"
$ws.Range("E41").Value = "'1.28%"
$ws.Range("D42").Value = "'0.004144"
$ws.Range("E42").Value = "'-33.59%"
$ws.Range("D43").Value = "'0.01503"
$ws.Range("E43").Value = "'25.46%"
$ws.Range("D44").Value = "'0.002221"
$ws.Range("E44").Value = "'16.89%"
$ws.Range("D45").Value = "'0.00005294"
$ws.Range("E45").Value = "'-3.10%"
$ws.Range("E48").Value = "'-46.99%"
